$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate new rows 74-99 with Tu Vi readings (column A = condition, column B = result)
$ws.Range("A74").Value = 'Sinh năm Giáp có Tử Vi đồng cung Thiên Phủ tại Dần'
$ws.Range("B74").Value = 'Sinh năm Giáp có Tử Vi đồng cung Thiên Phủ tại Dần'
$ws.Range("A75").Value = 'Tử Vi đồng cung Thiên Phủ tại Dần'
$ws.Range("B75").Value = 'Tử Vi đồng cung Thiên Phủ tại Dần'
$ws.Range("A76").Value = 'Tử Vi đồng cung Thiên Phủ tại Thân'
$ws.Range("B76").Value = 'Tử Vi đồng cung Thiên Phủ tại Thân'
$ws.Range("A77").Value = 'Sinh năm Giáp có Tử Vi đồng cung Thiên Phủ tại Dần'
$ws.Range("B77").Value = 'Sinh năm Giáp có Tử Vi đồng cung Thiên Phủ tại Dần'
$ws.Range("A78").Value = 'Sinh năm Giáp có Tử Vi đồng cung Thiên Phủ tại Thân'
$ws.Range("B78").Value = 'Sinh năm Giáp có Tử Vi đồng cung Thiên Phủ tại Thân'
$ws.Range("A79").Value = 'Tử Vi tọa thủ cung Mệnh và hội chiếu Thiên Phủ gặp Tả Phù, Hữu Bật'
$ws.Range("B79").Value = 'Được hưởng phú quý trọn đời.'
$ws.Range("A80").Value = 'Tử Vi tọa thủ cung Mệnh và hội chiếu Thiên Phủ có Kình Dương đồng cung'
$ws.Range("B80").Value = 'Kinh doanh buôn bán đại phú.'
$ws.Range("A81").Value = 'ử Phủ Vũ Tướng Tả Hữu Khoa Quyền Lộc Long Phượng'
$ws.Range("B81").Value = 'Hưởng đại phú đến cực độ, tuổi thọ gia tăng.'
$ws.Range("A82").Value = ' Tử Vi tọa thủ cung Mệnh và gặp Thiên Phủ, Vũ Khúc, Thiên Tướng, Tả Phù, Hữu Bật, Long Trì, Phượng Các, Hóa Khoa, Hóa Quyền, Hóa Lộc không gặp Kình Dương, Thiên Kiếp'
$ws.Range("B82").Value = ' Tử Vi tọa thủ cung Mệnh và gặp Thiên Phủ, Vũ Khúc, Thiên Tướng, Tả Phù, Hữu Bật, Long Trì, Phượng Các, Hóa Khoa, Hóa Quyền, Hóa Lộc không gặp Kình Dương, Thiên Kiếp'
$ws.Range("A83").Value = 'Sinh năm Kỷ có Tử Vi đồng cung Thiên Phủ tại Dần tại Mệnh gặp Hóa Quyền'
$ws.Range("B83").Value = 'Sinh năm Kỷ có Tử Vi đồng cung Thiên Phủ tại Dần tại Mệnh gặp Hóa Quyền'
$ws.Range("A84").Value = 'Tử Vi Thất Sát đồng cung tọa thủ cung Mệnh ở Tỵ'
$ws.Range("B84").Value = 'Tử Vi Thất Sát đồng cung tọa thủ cung Mệnh ở Tỵ'
$ws.Range("A85").Value = 'Tử Vi Thất Sát đồng cung tọa thủ cung Mệnh ở Hợi'
$ws.Range("B85").Value = 'Tử Vi Thất Sát đồng cung tọa thủ cung Mệnh ở Hợi'
$ws.Range("A86").Value = 'Tử Vi Thất Sát đồng cung tọa thủ cung Mệnh ở Hợi gặp Hóa Quyền'
$ws.Range("B86").Value = 'Tử Vi Thất Sát đồng cung tọa thủ cung Mệnh ở Hợi gặp Hóa Quyền'
$ws.Range("A87").Value = 'Tử Vi Thất Sát đồng cung tọa thủ cung Mệnh ở Tỵ gặp Hóa Quyền'
$ws.Range("B87").Value = 'Tử Vi Thất Sát đồng cung tọa thủ cung Mệnh ở Tỵ gặp Hóa Quyền'
$ws.Range("A88").Value = 'Tử Vi tọa thủ cung Mệnh và gặp Vũ Khúc, Phá Quân, Kình Dương, Đà La'
$ws.Range("B88").Value = 'Tử Vi tọa thủ cung Mệnh và gặp Vũ Khúc, Phá Quân, Kình Dương, Đà La'
$ws.Range("A89").Value = 'Tử Vi tọa thủ cung Mệnh ở Thìn có Phá Toái đồng cung'
$ws.Range("B89").Value = 'Tử Vi tọa thủ cung Mệnh ở Thìn có Phá Toái đồng cung'
$ws.Range("A90").Value = 'Tử Vi tọa thủ cung Mệnh ở Tuất có Phá Toái đồng cung'
$ws.Range("B90").Value = 'Tử Vi tọa thủ cung Mệnh ở Tuất có Phá Toái đồng cung'
$ws.Range("A91").Value = 'Tử Vi tọa thủ cung Mệnh ở Sửu có Phá Toái đồng cung'
$ws.Range("B91").Value = 'Tử Vi tọa thủ cung Mệnh ở Sửu có Phá Toái đồng cung'
$ws.Range("A92").Value = 'Tử Vi tọa thủ cung Mệnh ở Thìn có Phá Toái đồng cung'
$ws.Range("B92").Value = 'Tử Vi tọa thủ cung Mệnh ở Thìn có Phá Toái đồng cung'
$ws.Range("A93").Value = 'Tử Vi tọa thủ cung Mệnh có Phá Toái đồng cung với Địa Không, Địa Kiếp'
$ws.Range("B93").Value = 'Tử Vi tọa thủ cung Mệnh có Phá Toái đồng cung với Địa Không, Địa Kiếp'
$ws.Range("A94").Value = 'Tử Vi tọa thủ cung Mệnh có Đào Hoa, Hồng Loan, Địa Không, Địa Kiếp'
$ws.Range("B94").Value = 'Tử Vi tọa thủ cung Mệnh có Đào Hoa, Hồng Loan, Địa Không, Địa Kiếp'
$ws.Range("A95").Value = 'Tử Vi tọa thủ cung Mệnh và hội chiếu các sao Văn Khúc, Văn Xương'
$ws.Range("B95").Value = 'Tử Vi tọa thủ cung Mệnh và hội chiếu các sao Văn Khúc, Văn Xương'
$ws.Range("A96").Value = 'Tử Vi đồng cung Tả Phù Hữu Bật'
$ws.Range("B96").Value = 'Tử Vi đồng cung Tả Phù Hữu Bật'
$ws.Range("A97").Value = 'Tử Vi tọa thủ cung Mệnh ở Tý gặp Quyền, Lộc, Khoa'
$ws.Range("B97").Value = 'Tử Vi tọa thủ cung Mệnh ở Tý gặp Quyền, Lộc, Khoa'
$ws.Range("A98").Value = 'Tử Vi tọa thủ cung Mệnh ở Ngọ gặp Quyền, Lộc, Khoa'
$ws.Range("B98").Value = 'Tử Vi tọa thủ cung Mệnh ở Ngọ gặp Quyền, Lộc, Khoa'
$ws.Range("A99").Value = 'Tử Vi tọa thủ cung Mệnh gặp Quyền, Lộc, Khoa, Kình, Đà'
$ws.Range("B99").Value = 'Tử Vi tọa thủ cung Mệnh gặp Quyền, Lộc, Khoa, Kình, Đà'

# Apply the same yellow highlight fill used by the existing condition/result rows (style index 1)
$ws.Range("A74:B99").Interior.Color = 65535

# Restore view state: scroll position + active selection, matching the saved workbook state
$ws.Range("B103").Select()
try {
    $excel.ActiveWindow.ScrollRow = 79
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}
